# Apply the edits described by the diff:
# 1. Q62: 2 -> 0
# 2. Q71: 2 -> 0
# 3. O1146: 0 -> 2
# 4. R1148: (empty) -> 0
# 5. R1149: (empty) -> 0
# 6. Append new weekly rows 1150-1173 (columns A-Q; R left blank like the
#    other newly-appended rows in the source diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- single-cell fixes -----------------------------------------------
$ws.Cells.Item(62, 17).Value = 0
$ws.Cells.Item(71, 17).Value = 0
$ws.Cells.Item(1146, 15).Value = 2
$ws.Cells.Item(1148, 18).Value = 0
$ws.Cells.Item(1149, 18).Value = 0

# --- append new rows 1150-1173 ----------------------------------------
# Each inner array is: RowNumber, A(Datetime), B(Open), C(High), D(Low),
# E(Close), F(Adj Close), G(Volume), H(Year), I(Month), J(Day), K(Hour),
# L(Minute), M(Second), N(Week), O(isPivot), P(two_line_structure),
# Q(detect_structure)
$newRows = @(
    @(1150, 45474, 6180, 6342.85009765625, 6077, 6328.5498046875, 6318.83740234375, 1571523, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(1151, 45481, 6380, 6419.9501953125, 6255.89990234375, 6353.89990234375, 6344.1484375, 1355665, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(1152, 45488, 6340.5, 6539.7001953125, 6340.5, 6381.9501953125, 6372.15576171875, 1469826, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(1153, 45495, 6381.9501953125, 6679.9501953125, 6324.7998046875, 6664.5498046875, 6654.32177734375, 1661602, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(1154, 45502, 6777, 6795.85009765625, 6584.10009765625, 6710.14990234375, 6699.85205078125, 1934199, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(1155, 45509, 6689, 6768, 6500, 6587.85009765625, 6577.73974609375, 2226117, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(1156, 45516, 6550, 6865, 6420.9501953125, 6686.89990234375, 6676.6376953125, 2964678, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(1157, 45523, 6750, 6927.39990234375, 6565.0498046875, 6860.7001953125, 6860.7001953125, 1649010, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(1158, 45530, 6830, 6990, 6608.7998046875, 6927.5, 6927.5, 2092215, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(1159, 45537, 6989, 7044, 6790.0498046875, 6906.7998046875, 6906.7998046875, 1572837, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(1160, 45544, 6906.7998046875, 7099.9501953125, 6820.4501953125, 7019.14990234375, 7019.14990234375, 1951145, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(1161, 45551, 7025, 7152.4501953125, 6960.75, 7082.64990234375, 7082.64990234375, 1644934, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(1162, 45558, 7090.39990234375, 7279.89990234375, 7074.14990234375, 7234, 7234, 2002704, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(1163, 45565, 7160.5, 7316.9501953125, 6727.35009765625, 6773.60009765625, 6773.60009765625, 1838205, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(1164, 45572, 6801, 7070, 6743.14990234375, 6983.39990234375, 6983.39990234375, 1311693, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(1165, 45579, 6976.85009765625, 7150.35009765625, 6950, 6982.89990234375, 6982.89990234375, 1134816, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(1166, 45586, 6972.25, 7037.25, 6854.75, 6944.85009765625, 6944.85009765625, 766979, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0),
    @(1167, 45593, 6923, 7050, 6794.25, 7031.9501953125, 7031.9501953125, 869194, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0),
    @(1168, 45600, 7090, 7545, 6814.25, 7421.39990234375, 7421.39990234375, 4514544, 2024, 11, 4, 0, 0, 0, 45, 1, 0, 0),
    @(1169, 45607, 7400, 7454.14990234375, 6827.0498046875, 6860.64990234375, 6860.64990234375, 1312226, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(1170, 45614, 6874.4501953125, 6950.2998046875, 6594.14990234375, 6935.10009765625, 6935.10009765625, 1315931, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(1171, 45621, 6980, 7147, 6812, 6828.89990234375, 6828.89990234375, 1868745, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0),
    @(1172, 45628, 6858, 7334.9501953125, 6821.10009765625, 7233.2998046875, 7233.2998046875, 2235490, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0),
    @(1173, 45635, 7162.75, 7364, 7162.75, 7259.4501953125, 7259.4501953125, 1639902, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0)
)

$dateFormat = $ws.Cells.Item(1149, 1).NumberFormat

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
}
